# Weekly update: rows 5-19 shift down by one (oldest entry removed from row 5,
# a brand-new entry is written into row 5, and the entry that used to be in
# row 19 is appended as the new row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = 45054
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 14000
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 778

# Row 6
$ws.Range("D6").Value = 44328
$ws.Range("M6").Value = 250

# Row 7
$ws.Range("D7").Value = 44319
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 120

# Row 8
$ws.Range("D8").Value = 44714
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("Q8").Value = "`$/caja 18 kilos granel"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 1111
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44291
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("Q9").Value = "`$/caja 15 kilos granel"
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 800
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("D10").Value = 44691
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 17000
$ws.Range("P10").Value = 17000
$ws.Range("S10").Value = 944

# Row 11
$ws.Range("D11").Value = 44342

# Row 12
$ws.Range("D12").Value = 44316
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 300

# Row 13
$ws.Range("D13").Value = 44340
$ws.Range("M13").Value = 230
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("S13").Value = 1111

# Row 14
$ws.Range("D14").Value = 44354
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 18000
$ws.Range("S14").Value = 1000

# Row 15
$ws.Range("D15").Value = 44326
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("S15").Value = 1111

# Row 16
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("S16").Value = 1000

# Row 17
$ws.Range("D17").Value = 44358
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 17000
$ws.Range("P17").Value = 17000
$ws.Range("S17").Value = 944

# Row 18
$ws.Range("D18").Value = 44355
$ws.Range("L18").Value = "Especial"
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 18000
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 18000
$ws.Range("Q18").Value = "`$/caja 18 kilos granel"
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44680
$ws.Range("M19").Value = 200
$ws.Range("R19").Value = "Provincia de Limarí"

# New row 20 (appended)
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Macroferia Regional de Talca"
$ws.Range("C20").Value = "Maule"
$ws.Range("D20").Value = 44299
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100104
$ws.Range("H20").Value = "Frutos de pepita"
$ws.Range("I20").Value = 100104001
$ws.Range("J20").Value = "Granada"
$ws.Range("K20").Value = "Wonderfull"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 15000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 15000
$ws.Range("Q20").Value = "`$/caja 15 kilos granel"
$ws.Range("R20").Value = "Provincia de Curicó"
$ws.Range("S20").Value = 1000
$ws.Range("T20").Value = 15
